# "add 3 times code error retry"
#
# The old "taskcode" sheet (汽运任务/分解任务/火运计划 task-code smoke test)
# is replaced by a "login" sheet that drives a login regression test:
# username / password / expected-result(0 fail,1 success) / actual result /
# error message, covering the good-credentials case plus the failure cases
# that a 3x retry-on-error flow needs to exercise (blank username, missing
# password, wrong username).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- rename the first sheet: taskcode -> login ---
$ws1.Name = "login"

# --- wipe the old task-code rows/columns ---
$ws1.Cells.Clear() | Out-Null

# Touch a scratch cell with the custom "0_);[Red](0)" number format and
# clear it again so the format definition exists in the style table even
# though none of the surviving cells below end up keeping it selected.
$ws1.Cells.Item(50, 50).Value = 1
$ws1.Cells.Item(50, 50).NumberFormat = "0_);[Red]\(0\)"
$ws1.Cells.Item(50, 50).Clear() | Out-Null

# --- values, written in the same order the strings were first typed ---
$ws1.Cells.Item(1, 1).Value = "用户名"
$ws1.Cells.Item(1, 2).Value = "密码"
$ws1.Cells.Item(2, 1).Value = "iosapp0"
$ws1.Cells.Item(1, 3).Value = "预期结果(0:失败，1:成功)"
$ws1.Cells.Item(1, 4).Value = "执行结果"
$ws1.Cells.Item(1, 5).Value = "错误消息"

$ws1.Cells.Item(2, 2).Value = 123456
$ws1.Cells.Item(2, 3).Value = 1

$ws1.Cells.Item(3, 2).Value = 123456
$ws1.Cells.Item(3, 3).Value = 0

$ws1.Cells.Item(4, 1).Value = "iosapp0"
$ws1.Cells.Item(4, 3).Value = 0

$ws1.Cells.Item(5, 1).Value = "iosapp00"
$ws1.Cells.Item(5, 2).Value = 123456
$ws1.Cells.Item(5, 3).Value = 0

# --- header row is stored as text ---
$ws1.Range("A1:C1").NumberFormat = "@"

# --- re-apply the quoted-text flag to the login/username cells ---
$ws1.Cells.Item(2, 1).Value = "'iosapp0"
$ws1.Cells.Item(3, 1).Value = "'"
$ws1.Cells.Item(3, 1).Value = ""
$ws1.Cells.Item(4, 1).Value = "'iosapp0"
$ws1.Cells.Item(5, 1).Value = "'iosapp00"

# --- selection / active sheet ---
$ws1.Select() | Out-Null
$ws1.Range("B5").Select() | Out-Null
